# Regenerate the s_vals data (filtering save games), updating the
# per-row stat columns (TB, d2S, K, IP) and the recomputed "sum" column.
# The "Win" column (F) and date labels (A) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    3  = @{ B = 0.6753301551942219; C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732; G = 2.997429241610044 }
    4  = @{ B = 0.04763786555579896;C = 0.3127903958511391;D = 0.1575252929769615; E = 0.496779210170732; G = 1.014732764554632 }
    5  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    6  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    7  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    9  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    10 = @{ B = 0.01514828764759746;C = 0.04240448674262143;D = 0.8054896365839992;E = 0.496779210170732; G = 1.35982162114495 }
    11 = @{ B = 0.127881588408715;  C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732; G = 1.094976487407548 }
    12 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
